$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.077.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +5.67%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.718.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.57%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.11%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'332.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.55%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +0.08%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.3686"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.20%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'49.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.38%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.3346"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +2.42%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'1.185"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.52%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07470"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +5.71%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.31%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'6.292"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +4.96%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'20.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.60%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'6.922"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.46%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'1.715.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.43%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.00001077"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.85%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.06634"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.05%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'82.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +3.93%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.15%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  +4.02%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'6.087"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.59%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'13.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.96%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'26.012.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +5.54%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.476"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.16%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.460"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.61%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'150.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.53%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'19.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +3.19%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'1.323"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +8.86%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'1.905.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.46%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'129.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.14%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'4.103"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.67%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'5.922"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.31%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'0.08507"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.33%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'1.718"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.83%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'12.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +4.35%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'5.342"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.41%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06221"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +3.01%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02296"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.59%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2132"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.58%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'8.528"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +3.80%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'1.244"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.82%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'14.43"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +13.15%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "'0.6168"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.87%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.12%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'3.835"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.65%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.5886"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +4.54%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'127.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.65%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'2.018"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.02%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.07280"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +4.48%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'76.99"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.99%  "
$ws.Range("E51").Style = "Normal"
